$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Pete"
$ws.Range("C10").Value = "Sampras"

$ws.Range("C10").Select()
